$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 110
$ws1.Range("F3").Value = 7424
$ws1.Range("F4").Value = 276
$ws1.Range("F6").Value = 3969
$ws1.Range("F7").Value = 318
$ws1.Range("F8").Value = 561
$ws1.Range("F10").Value = 634
$ws1.Range("F11").Value = 118

# Sheet "演出" (shows)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 2

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 110
$ws4.Range("F4").Value = 7424
$ws4.Range("F6").Value = 276
$ws4.Range("F8").Value = 3969
$ws4.Range("F9").Value = 318
$ws4.Range("F10").Value = 561
$ws4.Range("F12").Value = 634
$ws4.Range("F13").Value = 2
$ws4.Range("F14").Value = 118

$wb.Save()
